# Append a new paragraph describing the cinnamon-apple pie at the very end
# of the document, right after the last (empty) paragraph and before the
# section properties, matching the author's added content.

$d = $word.ActiveDocument

# Locate the very end of the document content (after the last paragraph mark)
$endPos = $d.Content.End
$endRange = $d.Range($endPos, $endPos)

# Insert a brand-new paragraph after the current last (empty) paragraph.
# It inherits the same paragraph-level formatting (Arial font, color
# 222222, white shading) as the paragraph it follows.
$endRange.InsertParagraphAfter()

# Grab the newly created (now last) paragraph and place our text inside it.
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newRange = $newPara.Range
$newRange.Collapse(1)

$newRange.InsertAfter("Our cinnamon-apple pie is flavored with white and brown sugar, loads of cinnamon and maple syrup; its crust is golden brown and the pie filling bubbles through the lattice.")
